$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure price column (D) values are stored as text, preserving formats like trailing zeros
# and avoiding locale-based numeric parsing of values such as "63.760.99".
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = "63.760.99"
$ws.Range("E2").Value = "  +1.05%  "

$ws.Range("D3").Value = "2.611.23"
$ws.Range("E3").Value = "  +0.00%  "

$ws.Range("E4").Value = "  +0.06%  "

$ws.Range("D5").Value = "596.24"
$ws.Range("E5").Value = "  -1.57%  "

$ws.Range("D6").Value = "151.28"
$ws.Range("E6").Value = "  +4.08%  "

$ws.Range("E7").Value = "  +0.05%  "

$ws.Range("D8").Value = "0.588"
$ws.Range("E8").Value = "  +0.76%  "

$ws.Range("E9").Value = "  +1.83%  "

$ws.Range("D10").Value = "5.68"
$ws.Range("E10").Value = "  +3.38%  "

$ws.Range("D11").Value = "0.385"
$ws.Range("E11").Value = "  +3.70%  "

$ws.Range("E12").Value = "  -0.77%  "

$ws.Range("D13").Value = "27.95"
$ws.Range("E13").Value = "  +3.10%  "

$ws.Range("D14").Value = "3.087.06"
$ws.Range("E14").Value = "  +0.22%  "

$ws.Range("D15").Value = "63.631.53"
$ws.Range("E15").Value = "  +1.08%  "

$ws.Range("E16").Value = "  +6.07%  "

$ws.Range("D17").Value = "2.627.34"
$ws.Range("E17").Value = "  +1.43%  "

$ws.Range("D18").Value = "12.43"
$ws.Range("E18").Value = "  +8.30%  "

$ws.Range("D19").Value = "4.71"
$ws.Range("E19").Value = "  +4.77%  "

$ws.Range("D20").Value = "349.25"
$ws.Range("E20").Value = "  +2.25%  "

$ws.Range("D21").Value = "6.88"
$ws.Range("E21").Value = "  +0.50%  "

$ws.Range("D22").Value = "0.998"
$ws.Range("E22").Value = "  -0.25%  "

$ws.Range("D23").Value = "67.50"
$ws.Range("E23").Value = "  +2.13%  "

$ws.Range("D24").Value = "1.71"
$ws.Range("E24").Value = "  +7.92%  "

$ws.Range("D25").Value = "9.38"
$ws.Range("E25").Value = "  +4.34%  "

$ws.Range("D26").Value = "1.70"
$ws.Range("E26").Value = "  +1.02%  "

$ws.Range("D27").Value = "560.34"
$ws.Range("E27").Value = "  +3.22%  "

$ws.Range("D28").Value = "8.03"
$ws.Range("E28").Value = "  +2.97%  "

$ws.Range("D29").Value = "0.162"
$ws.Range("E29").Value = "  +0.33%  "

$ws.Range("E30").Value = "  -0.01%  "

$ws.Range("D31").Value = "2.06"
$ws.Range("E31").Value = "  +1.45%  "

$ws.Range("D32").Value = "0.0₃0853"
$ws.Range("E32").Value = "  +1.60%  "

$ws.Range("D33").Value = "1.75"
$ws.Range("E33").Value = "  +0.33%  "

$ws.Range("D34").Value = "5.32"
$ws.Range("E34").Value = "  +2.19%  "

$ws.Range("D35").Value = "166.66"
$ws.Range("E35").Value = "  -0.85%  "

$ws.Range("D36").Value = "0.416"
$ws.Range("E36").Value = "  +3.71%  "

$ws.Range("D37").Value = "0.999"
$ws.Range("E37").Value = "  -0.14%  "

$ws.Range("D38").Value = "19.63"
$ws.Range("E38").Value = "  +3.87%  "

$ws.Range("D39").Value = "1.95"
$ws.Range("E39").Value = "  +0.82%  "

$ws.Range("E40").Value = "  +0.02%  "

$ws.Range("D41").Value = "166.89"
$ws.Range("E41").Value = "  +1.18%  "

$ws.Range("D42").Value = "39.70"
$ws.Range("E42").Value = "  +0.18%  "

$ws.Range("D43").Value = "3.97"
$ws.Range("E43").Value = "  +6.09%  "

$ws.Range("D44").Value = "0.0588"
$ws.Range("E44").Value = "  +4.67%  "

$ws.Range("E45").Value = "  +1.78%  "

$ws.Range("D46").Value = "0.632"
$ws.Range("E46").Value = "  +1.42%  "

$ws.Range("B47").Value = "VeChain"
$ws.Range("C47").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D47").Value = "0.0252"
$ws.Range("E47").Value = "  +4.02%  "

$ws.Range("B48").Value = "dogwifhat"
$ws.Range("C48").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D48").Value = "2.04"
$ws.Range("E48").Value = "  +6.47%  "

$ws.Range("D49").Value = "0.0969"
$ws.Range("E49").Value = "  +1.46%  "

$ws.Range("D50").Value = "19.27"
$ws.Range("E50").Value = "  +3.42%  "

$ws.Range("D51").Value = "0.0₆0236"
$ws.Range("E51").Value = "  +20.88%  "
